$d = $word.ActiveDocument

# Locate the sentence we need to extend: "This is a Microsoft word document."
$findRange = $d.Content.Duplicate
$found = $findRange.Find.Execute("This is a Microsoft word document.", $true, $false, $false, `
                                  $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find target sentence"
}

$origText = $findRange.Text
$origTextEsc = $origText -replace '&', '&amp;' -replace '<', '&lt;' -replace '>', '&gt;'

# Build a small OOXML package that replaces the matched range's content
# with the original run PLUS three new runs: " (", "Changed main", ")".
# Because the replaced range stops right before the paragraph mark, the
# surrounding paragraph (with all of its original attributes) is left
# untouched -- only the run-level content inside it is substituted.
$bodyXml = '<w:p>' + `
    '<w:r><w:t>' + $origTextEsc + '</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> (</w:t></w:r>' + `
    '<w:r><w:t>Changed main</w:t></w:r>' + `
    '<w:r><w:t>)</w:t></w:r>' + `
    '</w:p>'

$pkgXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
    '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
    '<w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$findRange.InsertXML($pkgXml)
